{"js": "// Replace each two-digit multiplication equation's text with the new value.\n// Old values are unique in the document, so a direct search+replace per pair is safe and order-independent.\nconst replacements = [\n  [\"80\u00d780=6400\", \"49\u00d764=3136\"],\n  [\"72\u00d750=3600\", \"23\u00d738=874\"],\n  [\"32\u00d780=2560\", \"40\u00d797=3880\"],\n  [\"16\u00d741=656\", \"58\u00d786=4988\"],\n  [\"63\u00d739=2457\", \"75\u00d736=2700\"],\n  [\"25\u00d724=600\", \"64\u00d747=3008\"],\n  [\"50\u00d775=3750\", \"18\u00d714=252\"],\n  [\"46\u00d779=3634\", \"89\u00d722=1958\"],\n  [\"79\u00d756=4424\", \"54\u00d743=2322\"],\n  [\"97\u00d777=7469\", \"89\u00d745=4005\"],\n  [\"19\u00d765=1235\", \"31\u00d724=744\"],\n  [\"28\u00d720=560\", \"33\u00d787=2871\"],\n  [\"63\u00d781=5103\", \"71\u00d781=5751\"],\n  [\"20\u00d772=1440\", \"36\u00d799=3564\"],\n  [\"84\u00d751=4284\", \"49\u00d770=3430\"],\n  [\"41\u00d795=3895\", \"32\u00d720=640\"],\n  [\"38\u00d761=2318\", \"62\u00d731=1922\"],\n  [\"38\u00d748=1824\", \"63\u00d762=3906\"],\n  [\"74\u00d782=6068\", \"12\u00d790=1080\"],\n  [\"90\u00d784=7560\", \"25\u00d773=1825\"],\n  [\"78\u00d736=2808\", \"31\u00d782=2542\"],\n  [\"11\u00d789=979\", \"98\u00d721=2058\"],\n  [\"96\u00d743=4128\", \"66\u00d782=5412\"],\n  [\"77\u00d769=5313\", \"55\u00d756=3080\"],\n  [\"28\u00d781=2268\", \"31\u00d755=1705\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Update each two-digit multiplication equation to the new value via Find & Replace.\n# Old equation strings are unique in the document, so MatchWholeWord + ReplaceAll is safe per pair.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"80\u00d780=6400\", \"49\u00d764=3136\"),\n    @(\"72\u00d750=3600\", \"23\u00d738=874\"),\n    @(\"32\u00d780=2560\", \"40\u00d797=3880\"),\n    @(\"16\u00d741=656\", \"58\u00d786=4988\"),\n    @(\"63\u00d739=2457\", \"75\u00d736=2700\"),\n    @(\"25\u00d724=600\", \"64\u00d747=3008\"),\n    @(\"50\u00d775=3750\", \"18\u00d714=252\"),\n    @(\"46\u00d779=3634\", \"89\u00d722=1958\"),\n    @(\"79\u00d756=4424\", \"54\u00d743=2322\"),\n    @(\"97\u00d777=7469\", \"89\u00d745=4005\"),\n    @(\"19\u00d765=1235\", \"31\u00d724=744\"),\n    @(\"28\u00d720=560\", \"33\u00d787=2871\"),\n    @(\"63\u00d781=5103\", \"71\u00d781=5751\"),\n    @(\"20\u00d772=1440\", \"36\u00d799=3564\"),\n    @(\"84\u00d751=4284\", \"49\u00d770=3430\"),\n    @(\"41\u00d795=3895\", \"32\u00d720=640\"),\n    @(\"38\u00d761=2318\", \"62\u00d731=1922\"),\n    @(\"38\u00d748=1824\", \"63\u00d762=3906\"),\n    @(\"74\u00d782=6068\", \"12\u00d790=1080\"),\n    @(\"90\u00d784=7560\", \"25\u00d773=1825\"),\n    @(\"78\u00d736=2808\", \"31\u00d782=2542\"),\n    @(\"11\u00d789=979\", \"98\u00d721=2058\"),\n    @(\"96\u00d743=4128\", \"66\u00d782=5412\"),\n    @(\"77\u00d769=5313\", \"55\u00d756=3080\"),\n    @(\"28\u00d781=2268\", \"31\u00d755=1705\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}"}
